# Trade #23 closed at 2026-02-17 15:19:08 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status rollups for the MarketMaking
# strategy and appends the newly-closed trade (#23) as row 24 of both
# the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - refresh aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.8    # Current Capital
$summary.Range("B4").Value = -0.2      # Total P&L $
$summary.Range("B5").Value = -0.17     # Total P&L %
$summary.Range("B6").Value = 23        # Total Trades
$summary.Range("B7").Value = 6         # Winning Trades
$summary.Range("B9").Value = 26.09     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - refresh the MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.8       # Capital
$status.Range("D4").Value = 23         # Trades
$status.Range("E4").Value = -0.2       # P&L $
$status.Range("F4").Value = -0.2       # P&L %
$status.Range("G4").Value = 26.09      # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade row (row 24) to a given trades sheet.
# ---------------------------------------------------------------------
function Add-TradeRow24($ws) {
    $ws.Range("A24").Value = 23

    # "2026-02-17" must stay a text value, not get auto-converted to a
    # date serial number - format the cell as Text first, then put the
    # style back to Normal (matches the other rows, which carry no
    # explicit cell style).
    $ws.Range("B24").NumberFormat = "@"
    $ws.Range("B24").Value = "2026-02-17"
    $ws.Range("B24").Style = "Normal"

    $ws.Range("C24").Value = "15:19:01"
    $ws.Range("D24").Value = "MarketMaking"
    $ws.Range("E24").Value = "UP"
    $ws.Range("F24").Value = 0.83
    $ws.Range("G24").Value = 0.87
    $ws.Range("H24").Value = "CLOSED"
    $ws.Range("I24").Value = 4.8193
    $ws.Range("J24").Value = 0.04
    $ws.Range("K24").Value = 99.8
    $ws.Range("L24").Value = 0
    $ws.Range("M24").Value = 0
    $ws.Range("N24").Value = 0.6
    $ws.Range("O24").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P24").Value = "early_exit"
    $ws.Range("Q24").Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow24 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow24 $marketMaking
